# Added function to pass keySequence
# - Select the full header row on the "hub" sheet (no longer the active tab)
# - Add a new "system" worksheet after "testdata"
# - Copy the "hub" header row (with its style) into "system"
# - Add a data row with objectID "Title" / name_nl "systeem"
# - Make "system" the active sheet/tab with A3 selected

$wb = $excel.ActiveWorkbook

$hub = $wb.Worksheets.Item(1)
$hub.Range("A1:XFD1").Select() | Out-Null

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$system = $wb.Worksheets.Add($null, $lastSheet)
$system.Name = "system"

# Copy header row (values + style) from hub into system
$hub.Range("A1:M1").Copy($system.Range("A1:M1"))

$system.Range("A2").Value = "Title"
$system.Range("B2").Value = "systeem"

$system.Range("A3").Select() | Out-Null
$system.Activate() | Out-Null
